$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ColA($cellRef, $text) {
  $r = $ws.Range($cellRef)
  $r.Value = $text
  $r.Font.Bold = $true
  $r.VerticalAlignment = -4160
}

function Set-ColB($cellRef, $text) {
  $r = $ws.Range($cellRef)
  $r.Value = $text
  $r.Font.Bold = $false
  $r.WrapText = $true
  $r.VerticalAlignment = -4160
}

function Set-ColC($cellRef, $text) {
  $r = $ws.Range($cellRef)
  $r.Value = $text
  $r.Font.Color = 255
  $r.WrapText = $true
  $r.VerticalAlignment = -4160
}

# Row 10: Objetivos: -> new Portuguese objectives text
Set-ColB "B10" '1. Aproximar o aluno do seu aparelho vocal, ao nível da expressão falada ou cantada, em sua expressão individual ou coletiva (coral). 2. Propiciar ao aluno o contato com o repertório coral e seu papel no desenvolvimento da linguagem musical. 3. Propiciar ao aluno a oportunidade de fazer música vocal em conjunto, com o aprendizado técnico de parâmetros como: afinação, precisão, equilíbrio, fraseado etc.'
Set-ColC "C10" '1. Aproximar o aluno do seu aparelho vocal, ao nível da expressão falada ou cantada, em sua expressão individual ou coletiva (coral). 2. Propiciar ao aluno o contato com o repertório coral e seu papel no desenvolvimento da linguagem musical. 3. Propiciar ao aluno a oportunidade de fazer música vocal em conjunto, com o aprendizado técnico de parâmetros como: afinação, precisão, equilíbrio, fraseado etc.'

# Row 13: clear A13 (no label) and set B13/C13 to the teacher name
$ws.Range("A13").Clear()
Set-ColB "B13" '8452037 - Elisabeth Pinheiro da Silva Kondracki de Alcantara'
Set-ColC "C13" '8452037 - Elisabeth Pinheiro da Silva Kondracki de Alcantara'
$ws.Rows.Item(13).EntireRow.AutoFit()

# Row 14: Programa resumido: -> new Portuguese summary text
Set-ColA "A14" 'Programa resumido:'
Set-ColB "B14" 'Retomada e aperfeiçoamento dos repertórios e pressupostos técnicos da disciplina precedente de Canto Coral. Classificação Vocal. Respiração para o canto. Colocação da emissão em "Bocca Chiusa". O canto coral em uníssono.O canto coral em cânone. O canto coral em outras formações polifônicas. Leitura coral.'
Set-ColC "C14" 'Retomada e aperfeiçoamento dos repertórios e pressupostos técnicos da disciplina precedente de Canto Coral. Classificação Vocal. Respiração para o canto. Colocação da emissão em "Bocca Chiusa". O canto coral em uníssono.O canto coral em cânone. O canto coral em outras formações polifônicas. Leitura coral.'
$ws.Rows.Item(14).RowHeight = 60

# Row 15: Short syllabus: -> (shifted up from old row 14)
Set-ColA "A15" 'Short syllabus:'
Set-ColB "B15" 'Resumption and improvement of repertoires and technical assumptions of previous discipline Choir Singing.Vocal classification. Breathing into the singing. Placing the emission in Bocca Chiusa. Choral singing in unison. Choral singing in canon. Choral singing in other polyphonic formations. Coral reading.'
Set-ColC "C15" 'Resumption and improvement of repertoires and technical assumptions of previous discipline Choir Singing.Vocal classification. Breathing into the singing. Placing the emission in Bocca Chiusa. Choral singing in unison. Choral singing in canon. Choral singing in other polyphonic formations. Coral reading.'
$ws.Rows.Item(15).RowHeight = 60

# Row 16: Programa: -> new full Portuguese programa text
Set-ColA "A16" 'Programa:'
Set-ColB "B16" 'Retomada e aperfeiçoamento dos repertórios e pressupostos técnicos da disciplina precedente de Canto Coral. Classificação Vocal. Respiração para o canto - exercícios para localização da respiração baixa e média. Apoio e coluna de ar. Colocação da emissão em  "Bocca Chiusa": relaxamento da mandíbula, posicionamento da língua, suspensão do palato mole, exercícios de percepção do local onde a voz está se colocando, conexão do apoio e emissão vocal, passagem da "Bocca Chiusa" para vogais e outros sons nasais e guturais. O canto coral em uníssono: afinação, uniformidade tímbrica, precisão rítmica. O canto coral em cânone. O canto coral em outras formações polifônicas. Leitura coral: testagem da escuta harmônica e afinação à 1ª vista a várias vozes, memorização, fundamentos de teoria musical. Montagem e aperfeiçoamento de peças musicais - promovendo a aplicação das técnicas aprendidas. Conexão entre diafragma e emissão vocal.'
Set-ColC "C16" 'Retomada e aperfeiçoamento dos repertórios e pressupostos técnicos da disciplina precedente de Canto Coral. Classificação Vocal. Respiração para o canto - exercícios para localização da respiração baixa e média. Apoio e coluna de ar. Colocação da emissão em  "Bocca Chiusa": relaxamento da mandíbula, posicionamento da língua, suspensão do palato mole, exercícios de percepção do local onde a voz está se colocando, conexão do apoio e emissão vocal, passagem da "Bocca Chiusa" para vogais e outros sons nasais e guturais. O canto coral em uníssono: afinação, uniformidade tímbrica, precisão rítmica. O canto coral em cânone. O canto coral em outras formações polifônicas. Leitura coral: testagem da escuta harmônica e afinação à 1ª vista a várias vozes, memorização, fundamentos de teoria musical. Montagem e aperfeiçoamento de peças musicais - promovendo a aplicação das técnicas aprendidas. Conexão entre diafragma e emissão vocal.'
$ws.Rows.Item(16).RowHeight = 120

# Row 17: Syllabus: -> gains B17/C17 (EN full syllabus)
Set-ColA "A17" 'Syllabus:'
Set-ColB "B17" 'Resumption and improvement of repertoires and technical assumptions of previous discipline Choir Singing. Vocal classification.Breathing into the singing - Exercises for locating low and average breathing. Support and air column. Placing the emission in Bocca Chiusa. - relaxation of the mandible, tongue position, the soft palate suspension, local perception exercises for where the voice is being put, support connection and vocal emission, passing the Bocca Chiusa for vowels and other nasal and guttural sounds. Choral singing in unison. - The choral singing in unison, tuning, timbre uniformity, rhythmic precision. Choral singing in canon. Choral singing in other polyphonic formations. Coral reading.- testing of harmonic relaying and listening to 1st view the various voices, memorization, music theory basics. Assembling and improvement of pieces - promoting the application of learnt techniques. Connection between diaphragm and vocal emission.'
Set-ColC "C17" 'Resumption and improvement of repertoires and technical assumptions of previous discipline Choir Singing. Vocal classification.Breathing into the singing - Exercises for locating low and average breathing. Support and air column. Placing the emission in Bocca Chiusa. - relaxation of the mandible, tongue position, the soft palate suspension, local perception exercises for where the voice is being put, support connection and vocal emission, passing the Bocca Chiusa for vowels and other nasal and guttural sounds. Choral singing in unison. - The choral singing in unison, tuning, timbre uniformity, rhythmic precision. Choral singing in canon. Choral singing in other polyphonic formations. Coral reading.- testing of harmonic relaying and listening to 1st view the various voices, memorization, music theory basics. Assembling and improvement of pieces - promoting the application of learnt techniques. Connection between diaphragm and vocal emission.'
$ws.Rows.Item(17).RowHeight = 120

# Row 18: Avaliacao: -> clear stray B18/C18 (label-only row)
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Rows.Item(18).EntireRow.AutoFit()

# Row 19: Metodo: -> new text (the "A cada semestre..." paragraph)
Set-ColB "B19" 'A cada semestre é proposto um programa com cerca de 8 (oito) peças, sendo duas ou três de semestres anteriores e, consequentemente, cinco ou seis inéditas a ser apresentado pelo CORAL da EEL-USP em performances públicas definidas durante o período letivo.'
Set-ColC "C19" 'A cada semestre é proposto um programa com cerca de 8 (oito) peças, sendo duas ou três de semestres anteriores e, consequentemente, cinco ou seis inéditas a ser apresentado pelo CORAL da EEL-USP em performances públicas definidas durante o período letivo.'
$ws.Rows.Item(19).RowHeight = 60

# Row 20: Criterio: -> new text (the "Sendo uma atividade..." paragraph)
Set-ColB "B20" 'Sendo uma atividade prática e de grupo, fica inviável a realização de provas ou outras formas similares de avaliação. Esta se dará no dia a dia do aluno, levando em conta: assiduidade, pontualidade e material completo na pasta; participação construtiva em sala de aula e nas apresentações públicas - prontidão, envolvimento e seu real aproveitamento vocal e musical.'
Set-ColC "C20" 'Sendo uma atividade prática e de grupo, fica inviável a realização de provas ou outras formas similares de avaliação. Esta se dará no dia a dia do aluno, levando em conta: assiduidade, pontualidade e material completo na pasta; participação construtiva em sala de aula e nas apresentações públicas - prontidão, envolvimento e seu real aproveitamento vocal e musical.'
$ws.Rows.Item(20).RowHeight = 60

# Row 21: Norma de recuperacao: -> "nao tem"
Set-ColB "B21" 'não tem'
Set-ColC "C21" 'não tem'
$ws.Rows.Item(21).RowHeight = 60

# Row 22: Bibliografia: -> new bibliography text
Set-ColA "A22" 'Bibliografia:'
Set-ColB "B22" 'BEHLAU, Mara; REHDER, Maria Inês. Higiene vocal para o Canto Coral. Rio de Janeiro: Revinter, 1997.COELHO, Helena Wöhl. Técnica vocal para coros. 7ª Edição. São Leopoldo: Sinodal, 1994.CORBIN, Lynn Ann. Vocal pedagogy in the choral rehearsal: The selected concepts on choral tone quality, understanding of the singing process, and attitudes toward choir participation. 1982. 119 f. Tese (Doutorado em Filosofia) - Graduate School of The Ohio State University, Ohio State University, 1982.COSTA, Paulo Rubens Moraes. Diagnose em Canto Coral: parâmetros para análise e ferramentas para a avaliação. 2005. Dissertação (Mestrado em Musicologia) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2005.DRAHAN, Snizhana. Ouvir a voz: a percepção da produção vocal pelo regente coral - método eformação. 2007. 146 f. Dissertação (Mestrado em Musicologia) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2007.FERNANDES, Angelo José. O regente moderno e a construção da sonoridade coral: uma metodologia de preparo vocal para coros. 2009. 475 f. Tese (Doutorado em Música) - Instituto de Artes, Universidade Federal de Campinas, Campinas, 2009.HERR, Martha. Considerações para a classificação da voz do coralista. In: FERREIRA, LésliePiccolotto et al. Voz profissional: o profissional da voz. Carapicuíba: Pró-fono DepartamentoEditorial, 1995. p. 51-56.MILLER, Richard. The structure of singing. Boston, Massachusetts: Schirmer, 1986.RAMOS, Marco Antonio da Silva. Canto Coral: do repertório temático à construção do programa. 1988. 492f. Dissertação (Mestrado em Artes) - Departamento de Biblioteconomia e Documentação, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 1988.______. Memorial. 2011. Memorial (Professor titular) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2011.______. O ensino da regência coral. 2003. 107f. Tese (Livre Docência) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2003.VIDEIRA JR., Mário Rodrigues. Educação musical através do coro: A experiência no Projeto Comunicantus da ECA - USP. 2001. 65 f. Trabalho de Conclusão de Curso (Licenciatura em Música) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2001.'
Set-ColC "C22" 'BEHLAU, Mara; REHDER, Maria Inês. Higiene vocal para o Canto Coral. Rio de Janeiro: Revinter, 1997.COELHO, Helena Wöhl. Técnica vocal para coros. 7ª Edição. São Leopoldo: Sinodal, 1994.CORBIN, Lynn Ann. Vocal pedagogy in the choral rehearsal: The selected concepts on choral tone quality, understanding of the singing process, and attitudes toward choir participation. 1982. 119 f. Tese (Doutorado em Filosofia) - Graduate School of The Ohio State University, Ohio State University, 1982.COSTA, Paulo Rubens Moraes. Diagnose em Canto Coral: parâmetros para análise e ferramentas para a avaliação. 2005. Dissertação (Mestrado em Musicologia) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2005.DRAHAN, Snizhana. Ouvir a voz: a percepção da produção vocal pelo regente coral - método eformação. 2007. 146 f. Dissertação (Mestrado em Musicologia) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2007.FERNANDES, Angelo José. O regente moderno e a construção da sonoridade coral: uma metodologia de preparo vocal para coros. 2009. 475 f. Tese (Doutorado em Música) - Instituto de Artes, Universidade Federal de Campinas, Campinas, 2009.HERR, Martha. Considerações para a classificação da voz do coralista. In: FERREIRA, LésliePiccolotto et al. Voz profissional: o profissional da voz. Carapicuíba: Pró-fono DepartamentoEditorial, 1995. p. 51-56.MILLER, Richard. The structure of singing. Boston, Massachusetts: Schirmer, 1986.RAMOS, Marco Antonio da Silva. Canto Coral: do repertório temático à construção do programa. 1988. 492f. Dissertação (Mestrado em Artes) - Departamento de Biblioteconomia e Documentação, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 1988.______. Memorial. 2011. Memorial (Professor titular) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2011.______. O ensino da regência coral. 2003. 107f. Tese (Livre Docência) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2003.VIDEIRA JR., Mário Rodrigues. Educação musical através do coro: A experiência no Projeto Comunicantus da ECA - USP. 2001. 65 f. Trabalho de Conclusão de Curso (Licenciatura em Música) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2001.'
$ws.Rows.Item(22).RowHeight = 120

# Row 23: Requisitos: (moved up); clear old B23/C23
Set-ColA "A23" 'Requisitos:'
$ws.Range("B23").Clear()
$ws.Range("C23").Clear()
$ws.Rows.Item(23).EntireRow.AutoFit()

# Row 24 (new row): requirement course text
Set-ColB "B24" '8800010 -  Canto Coral II  (Requisito)
'
Set-ColC "C24" '8800010 -  Canto Coral II  (Requisito)
'
$ws.Rows.Item(24).RowHeight = 30

Write-Host "Edit complete"
